# Apply the "#5: insurance, claim, debt, investment done" edit.
#
# The 債務 (debt) sheet and 事業投資 (business investment) sheet had a
# header row (row 1) that was a stray duplicate of the first data row
# instead of real field names, and they were missing the common
# metadata columns (property_category, category, date, legislator_name,
# legislator_id, source_file, index) that the other sheets (car, land,
# building) already carry. This script fixes the header row and appends
# those metadata columns to every data row.

$wb = $excel.ActiveWorkbook

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------
# Sheet "債務" (debt) -> Worksheets index 4
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)

# --- Fix header row (row 1) to contain real field names -------
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"

# --- Add the new metadata header cells (H1:N1), formatted like G1 --
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Add the metadata values for each data row (rows 2-5) ------
$debtRows = 2,3,4,5
foreach ($r in $debtRows) {
    $ws.Range("G$r").Copy()
    $ws.Range("H$r`:N$r").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws.Range("H$r").Value = "normal"
    $ws.Range("I$r").Value = "debt"

    # "date" looks like an ISO date, so Excel would otherwise convert it
    # to a date serial number; force it to be stored as text instead,
    # then restore the plain (unformatted) look of the cell.
    $ws.Range("J$r").NumberFormat = "@"
    $ws.Range("J$r").Value = "2012-04-26"
    $ws.Range("G$r").Copy()
    $ws.Range("J$r").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws.Range("K$r").Value = "楊應雄"
    $ws.Range("L$r").Value = 1758
    $ws.Range("M$r").Value = "tmp248f1"
    $ws.Range("N$r").Value = $ws.Range("A$r").Value2
}

# ---------------------------------------------------------------
# Sheet "事業投資" (business investment) -> Worksheets index 5
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(5)

# --- Fix header row (row 1) to contain real field names -------
$ws2.Range("B1").Value = "owner"
$ws2.Range("C1").Value = "company"
$ws2.Range("D1").Value = "address"
$ws2.Range("E1").Value = "total"
$ws2.Range("F1").Value = "register_date"
$ws2.Range("G1").Value = "register_reason"

# --- Add the new metadata header cells (H1:N1), formatted like G1 --
$ws2.Range("G1").Copy()
$ws2.Range("H1:N1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws2.Range("H1").Value = "property_category"
$ws2.Range("I1").Value = "category"
$ws2.Range("J1").Value = "date"
$ws2.Range("K1").Value = "legislator_name"
$ws2.Range("L1").Value = "legislator_id"
$ws2.Range("M1").Value = "source_file"
$ws2.Range("N1").Value = "index"

# --- Add the metadata values for each data row (rows 2-4) ------
$investRows = 2,3,4
foreach ($r in $investRows) {
    $ws2.Range("G$r").Copy()
    $ws2.Range("H$r`:N$r").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws2.Range("H$r").Value = "normal"
    $ws2.Range("I$r").Value = "investment"

    # same ISO-date-as-text workaround as above
    $ws2.Range("J$r").NumberFormat = "@"
    $ws2.Range("J$r").Value = "2012-04-26"
    $ws2.Range("G$r").Copy()
    $ws2.Range("J$r").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws2.Range("K$r").Value = "楊應雄"
    $ws2.Range("L$r").Value = 1758
    $ws2.Range("M$r").Value = "tmp248f1"
    $ws2.Range("N$r").Value = $ws2.Range("A$r").Value2
}
